# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on Sheet1 to match the latest scraped cryptocurrency snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (matches the workbook's existing
# inline-string cells) even when the text looks like a number, and restore
# the cell's style afterwards so no incidental formatting change is introduced.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.224.03"
$ws.Range("E2").Value = "  -0.56%  "

Set-TextValue $ws.Range("D3") "1.588.41"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  -0.10%  "

Set-TextValue $ws.Range("D5") "211.79"
$ws.Range("E5").Value = "  +0.76%  "

Set-TextValue $ws.Range("D6") "0.504"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.31%  "

$ws.Range("E9").Value = "  -1.07%  "

$ws.Range("E10").Value = "  -1.63%  "

Set-TextValue $ws.Range("D11") "0.0845"
$ws.Range("E11").Value = "  -0.07%  "

Set-TextValue $ws.Range("D12") "1.811.67"
$ws.Range("E12").Value = "  -0.21%  "

Set-TextValue $ws.Range("D13") "1.600.05"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("E14").Value = "  -1.46%  "

Set-TextValue $ws.Range("D16") "63.82"
$ws.Range("E16").Value = "  -0.90%  "

Set-TextValue $ws.Range("D17") "26.233.10"

$ws.Range("E18").Value = "  -0.35%  "

Set-TextValue $ws.Range("D19") "7.44"
$ws.Range("E19").Value = "  -0.47%  "

Set-TextValue $ws.Range("D20") "213.98"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("E24").Value = "  -1.67%  "

Set-TextValue $ws.Range("D25") "144.61"
$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("E28").Value = "  -0.98%  "

Set-TextValue $ws.Range("D29") "15.10"
$ws.Range("E29").Value = "  -1.00%  "

Set-TextValue $ws.Range("D30") "0.0494"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("E32").Value = "  -1.05%  "

Set-TextValue $ws.Range("D33") "1.418.05"
$ws.Range("E33").Value = "  +8.34%  "

Set-TextValue $ws.Range("D34") "2.95"
$ws.Range("E34").Value = "  -1.54%  "

$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("E39").Value = "  +4.84%  "

Set-TextValue $ws.Range("D40") "0.821"
$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("E41").Value = "  -0.10%  "

Set-TextValue $ws.Range("D42") "0.936"
$ws.Range("E42").Value = "  -14.95%  "

Set-TextValue $ws.Range("D43") "0.764"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("E44").Value = "  -0.49%  "

Set-TextValue $ws.Range("D45") "1.723.47"
$ws.Range("E45").Value = "  -0.28%  "

Set-TextValue $ws.Range("D46") "61.15"
$ws.Range("E46").Value = "  -2.34%  "

Set-TextValue $ws.Range("D47") "85.87"
$ws.Range("E47").Value = "  -2.34%  "

$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("E49").Value = "  -0.31%  "

Set-TextValue $ws.Range("D50") "0.0967"
$ws.Range("E50").Value = "  -1.48%  "

Set-TextValue $ws.Range("D51") "0.999"
$ws.Range("E51").Value = "  -0.14%  "
